# Generate Report for Handback
# The handoff file "7a82f988-c0d5-4b60-b7af-c5251f112ba7.md" has moved from
# "Ready for handoff" to "Handed back: in sync with en-US" for both the
# zh-cn and de-de locales. Update the Overview sheet and each locale sheet
# accordingly, including clearing the stale "Error Detail" message and
# refreshing the "Latest Handback DateTime" timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ovWs = $wb.Worksheets.Item("Overview")
$ovWs.Range("E3").Value = "Handed back: in sync with en-US"
$ovWs.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zhWs = $wb.Worksheets.Item("zh-cn")
$zhWs.Range("C3").Value = "Handed back: in sync with en-US"
$zhWs.Range("K3").Value = "2016-08-22 02:58:25"
$zhWs.Range("P3").Value = ""
$zhWs.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet ---
$deWs = $wb.Worksheets.Item("de-de")
$deWs.Range("C3").Value = "Handed back: in sync with en-US"
$deWs.Range("K3").Value = "2016-08-22 02:58:32"
$deWs.Range("P3").Value = ""
$deWs.Columns.Item(16).ColumnWidth = 13.7470528738839
